$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("Data")
$wsCode = $wb.Worksheets.Item("Codebook")

# --- Add new codebook rows describing the Cardio / caff variables ---
$wsCode.Range("A5").Value = "Cardio"
$wsCode.Range("B5").Value = "minutes of cardio in one week"
$wsCode.Range("C5").Value = "whole integers "

$wsCode.Range("A6").Value = "caff"
$wsCode.Range("B6").Value = "number of cafienated beverages daily"
$wsCode.Range("C6").Value = "whole integers "

# --- Update the "Data" sheet header labels (row 1) to match the codebook ---
# "Weekly Cardio in Minutes" -> "Cardio"
$wsData.Range("D1").Value = "Cardio"
# "Daily Caffienated Beverages" -> "caff"
$wsData.Range("E1").Value = "caff"

# --- Update view/selection state: Codebook becomes the active sheet ---
$wsData.Range("E1").Select()

$wsCode.Activate()
$wsCode.Range("C15").Select()
